# 🔄 MAJ automatique BRVM via GitHub Actions
# Refresh des donnees de recommandations et du top YTD.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Feuille "Recommandations" : Titre, Jours en Hausse, Jours en Baisse,
# Variation Totale (%), Derniere Variation (%), Recommandation, Strategie
# -----------------------------------------------------------------
$recoData = @(
    ,@('BRVM - SERVICES PUBLICS', 0, 10, 4170.15, 105.86, '🟡 Observer', '➖ Neutre')
    ,@('NEI-CEDA CI', 0, 5, 3370, 690, '🟡 Observer', '➖ Neutre')
    ,@('AIR LIQUIDE CI', 0, 5, 3365, 625, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - AUTRES SECTEURS', 0, 5, 2905.57, 592.21, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - DISTRIBUTION', 0, 5, 2183.21, 460.26, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - AGRICULTURE', 0, 5, 1807.85, 364.03, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - TRANSPORT', 0, 5, 1801.53, 356.15, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - CONSOMMATION DISCRETIONNAIRE', 0, 5, 704.61, 152.85, '🟡 Observer', '➖ Neutre')
    ,@('BRVM-PRESTIGE', 0, 5, 684.11, 136.73, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - FINANCES', 0, 5, 678.36, 135.67, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - SERVICES FINANCIERS', 0, 5, 666.69, 133.34, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - INDUSTRIELS', 0, 5, 614.08, 126.36, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - ENERGIE', 0, 5, 540.43, 109.2, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - TELECOMMUNICATIONS', 0, 5, 475.57, 94.48, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - INDUSTRIE                  (**)', 0, 2, 414.76, 207.68, '🟡 Observer', '➖ Neutre')
    ,@('BRVM-PRINCIPAL                    (**)', 0, 2, 387.23, 193.83, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - CONSOMMATION DE BASE         (**)', 0, 2, 372.63, 187.07, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - INDUSTRIE                (**)', 0, 1, 206.88, 206.88, '🟡 Observer', '➖ Neutre')
    ,@('BRVM-PRINCIPAL                  (**)', 0, 1, 193.58, 193.58, '🟡 Observer', '➖ Neutre')
    ,@('BRVM - CONSOMMATION DE BASE              (**)', 0, 1, 187.55, 187.55, '🟡 Observer', '➖ Neutre')
    ,@('FILTISAC CI (FTSC)', 5, 0, 36.91, 7.4, '🟢 Achat', '✅ Renforcer')
    ,@('BERNABE CI (BNBC)', 3, 0, 20.46, 7.5, '🟢 Achat', '✅ Renforcer')
    ,@('CFAO MOTORS CI (CFAC)', 3, 1, 16.2, 7.35, '🟢 Achat', '✅ Renforcer')
    ,@('TRACTAFRIC MOTORS CI (PRSC)', 3, 1, 15.58, 7.39, '🟢 Achat', '✅ Renforcer')
    ,@('ORAGROUP TOGO (ORGT)', 2, 1, 7.34, -7.38, '🟡 Observer', '👀 À surveiller')
    ,@('VIVO ENERGY CI (SHEC)', 2, 0, 4.19, 2.92, '🟡 Observer', '➖ Neutre')
    ,@('ONATEL BF (ONTBF)', 1, 0, 2.61, 2.61, '🟡 Observer', '➖ Neutre')
    ,@('BICI CI (BICC)', 1, 0, 2.2, 2.2, '🟡 Observer', '➖ Neutre')
    ,@('NSIA BANQUE COTE D''IVOIRE (NSBC)', 1, 1, 0.49, -2.01, '🟡 Observer', '👀 À surveiller')
    ,@('SICABLE CI (CABC)', 1, 1, 0.25, 5.2, '🟡 Observer', '👀 À surveiller')
    ,@('TOTAL', 0, 5, 0, 0, '🟡 Observer', '➖ Neutre')
    ,@('UNIWAX CI (UNXC)', 1, 1, -0.25, 7.23, '🟡 Observer', '👀 À surveiller')
    ,@('SAPH CI (SPHC)', 1, 1, -0.32, 5.94, '🟡 Observer', '👀 À surveiller')
    ,@('BANK OF AFRICA ML (BOAM)', 0, 1, -1.15, -1.15, '🟡 Observer', '➖ Neutre')
    ,@('SOGB CI (SOGC)', 0, 1, -1.18, -1.18, '🟡 Observer', '➖ Neutre')
    ,@('AFRICA GLOBAL LOGISTICS CI (SDSC)', 0, 1, -1.68, -1.68, '🟡 Observer', '➖ Neutre')
    ,@('TOTALENERGIES MARKETING CI (TTLC)', 0, 1, -2.04, -2.04, '🟡 Observer', '➖ Neutre')
    ,@('ECOBANK COTE D''''IVOIRE (ECOC)', 0, 1, -2.23, -2.23, '🟡 Observer', '➖ Neutre')
    ,@('SOLIBRA CI (SLBC)', 0, 1, -2.75, -2.75, '🟡 Observer', '➖ Neutre')
    ,@('NEI-CEDA CI (NEIC)', 0, 1, -2.94, -2.94, '🟡 Observer', '➖ Neutre')
    ,@('ORANGE COTE D''IVOIRE (ORAC)', 0, 1, -3.34, -3.34, '🟡 Observer', '➖ Neutre')
    ,@('SETAO CI (STAC)', 0, 2, -4.81, -2.91, '🟡 Observer', '➖ Neutre')
    ,@('SUCRIVOIRE (SCRC)', 1, 3, -6.91, 6.78, '🔴 Vente', '⚠️ Risque de décrochage')
    ,@('AIR LIQUIDE CI (SIVC)', 0, 1, -7.35, -7.35, '🟡 Observer', '➖ Neutre')
    ,@('UNILEVER CI (UNLC)', 0, 1, -7.5, -7.5, '🟡 Observer', '➖ Neutre')
    ,@('LOTERIE NATIONALE DU BENIN (LNBB)', 0, 2, -14.92, -7.45, '🟡 Observer', '➖ Neutre')
)

$wsReco = $wb.Worksheets.Item("Recommandations")
$r = 2
foreach ($row in $recoData) {
    $wsReco.Cells.Item($r, 1).Value = $row[0]
    $wsReco.Cells.Item($r, 2).Value = $row[1]
    $wsReco.Cells.Item($r, 3).Value = $row[2]
    $wsReco.Cells.Item($r, 4).Value = $row[3]
    $wsReco.Cells.Item($r, 5).Value = $row[4]
    $wsReco.Cells.Item($r, 6).Value = $row[5]
    $wsReco.Cells.Item($r, 7).Value = $row[6]
    $r++
}

# -----------------------------------------------------------------
# Feuille "Top_YTD" : Titre, Progression YTD (%)
# -----------------------------------------------------------------
$topYtdData = @(
    ,@('BRVM - SERVICES PUBLICS', 141617521.45)
    ,@('NEI-CEDA CI', 2775468.15)
    ,@('AIR LIQUIDE CI', 2752742.69)
    ,@('BRVM - AUTRES SECTEURS', 1463340.99)
    ,@('BRVM - DISTRIBUTION', 444206.1)
    ,@('BRVM - AGRICULTURE', 209362.02)
    ,@('BRVM - TRANSPORT', 206527.66)
    ,@('BRVM - CONSOMMATION DISCRETIONNAIRE', 8001.1)
    ,@('BRVM-PRESTIGE', 7349.14)
    ,@('BRVM - FINANCES', 7170.07)
)

$wsTop = $wb.Worksheets.Item("Top_YTD")
$r = 2
foreach ($row in $topYtdData) {
    $wsTop.Cells.Item($r, 1).Value = $row[0]
    $wsTop.Cells.Item($r, 2).Value = $row[1]
    $r++
}

Write-Host "Recommandations et Top_YTD mis a jour."
